$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix style for C16/C63/C64 to match target (copy format from analogous cells)
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C62").Copy() | Out-Null
$ws.Range("C63").PasteSpecial(-4122) | Out-Null
$ws.Range("C64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new answer text across column C
$ws.Range("C2").Value = "Site4Now / Locaweb / Smarter"
$ws.Range("C3").Value = "( x  ) SAAS  ( x  ) PAAS ( x  ) IAAS"
$ws.Range("C4").Value = "(   ) Nuvem Pública   (   ) Nuvem Privada   ( x ) Hibrida"
$ws.Range("C5").Value = "Cloud Hosting / VPS"
$ws.Range("C6").Value = "Data Center com serviços Cloud, SO, Link e Banco de Dados"
$ws.Range("C8").Value = "TIER I"
$ws.Range("C9").Value = "América do Sul (São Paulo) / América do Norte (Los Angeles) / Europa (Amisterdã)"
$ws.Range("C11").Value = "Sim (Ambas)"
$ws.Range("C12").Value = "Sim (Ambas)"
$ws.Range("C14").Value = "Disponíveis equipes de Infra-estrutura e Segurança que monitoram o ambiente através de ferramentas, oferecendo dashboards com informações de disponibilidade e estado do ambiente."
$ws.Range("C15").Value = "Através de equipes especializadas nas áreas dentro do Data Center"
$ws.Range("C16").Value = "Através de equipes especializadas nas áreas dentro do Data Center"
$ws.Range("C17").Value = "Logs e Jobs com disparos de alertas"
$ws.Range("C18").Value = "Os servidores possuem ferramentas de firewalls, detecçao de vírus e malwares"
$ws.Range("C19").Value = "Através de relatórios periódicos para análise e alertas"
$ws.Range("C20").Value = "Equipamentos alocados dentro do data center"
$ws.Range("C22").Value = "Existem servidores virtualizados, e servidores separados para aplicação e banco de dados"
$ws.Range("C23").Value = "Utilização do Cloudflare"
$ws.Range("C24").Value = "Sim"
$ws.Range("C26").Value = "Sim"
$ws.Range("C27").Value = "São espelhados e armazenados em sites diferentes, em localidades distintas"
$ws.Range("C28").Value = "Até 24 horas"
$ws.Range("C29").Value = "Até 12 horas"
$ws.Range("C30").Value = "Plano atual realiza um backup full ao dia"
$ws.Range("C31").Value = "Sim em storages separados"
$ws.Range("C32").Value = "Backup full por um período de 7 dias"
$ws.Range("C33").Value = "Em Storages separados"
$ws.Range("C36").Value = "Roteadores e firewalls"
$ws.Range("C37").Value = "Aos adminstradores terão acesso full, restrito aos demais"
$ws.Range("C38").Value = "Existem redundância, não mantém servidores em DMZ"
$ws.Range("C39").Value = "Sim"
$ws.Range("C40").Value = "Através dos usuários cadastrados no próprio sistema operacional"
$ws.Range("C41").Value = "Sim"
$ws.Range("C42").Value = "Armazenadas em banco de dados, com criptografia em chave privada"
$ws.Range("C43").Value = "Sim, não trafega sem essa segurança"
$ws.Range("C44").Value = "Não"
$ws.Range("C45").Value = "Sim"
$ws.Range("C46").Value = "Não"
$ws.Range("C47").Value = "Implementada autenticação específica"
$ws.Range("C48").Value = "Exige apenas um mínimo de caracteres, possibilita troca, e envios da senha por e-mail"
$ws.Range("C49").Value = "Valida qtde mínima de caracteres"
$ws.Range("C50").Value = "Sim"
$ws.Range("C51").Value = "Dentro de regras no banco de dados"
$ws.Range("C52").Value = "Detecção, alarmes, extintores, dispositivos de desligamento automático, proteções, etc"
$ws.Range("C53").Value = "Dados são logicamente destruídos antes do descarte no ambiente apropriado"
$ws.Range("C54").Value = "Dados são logicamente destruídos antes do descarte no ambiente apropriado"
$ws.Range("C57").Value = "Não"
$ws.Range("C58").Value = "Certificados TLS / SSL"
$ws.Range("C59").Value = "Sim"
$ws.Range("C60").Value = "Utiliza dos padrões dos certificados"
$ws.Range("C62").Value = "Não exitem planos atuais no momento"
$ws.Range("C63").Value = "Não exitem planos atuais no momento"
$ws.Range("C64").Value = "Não exitem planos atuais no momento"
$ws.Range("C66").Value = "Profissionais com larga escala de experiência e responsabilidade no setor"
$ws.Range("C67").Value = "Através de Logs"
$ws.Range("C68").Value = "Monitoramento da equipe de segurança do data center"
$ws.Range("C70").Value = "Não exitem planos atuais no momento"
$ws.Range("C71").Value = "Sim"
$ws.Range("C73").Value = "Integrações via API Rest"
$ws.Range("C74").Value = "5 Mb"
$ws.Range("C75").Value = "20 Mb"
$ws.Range("C76").Value = "Além do navegador, o plugin assinador para a leitura do certificado digital"
$ws.Range("C77").Value = "Não"
$ws.Range("C79").Value = "Equipe de desenvolvimento disposto para melhorias e customizações"
$ws.Range("C80").Value = "Sim"
$ws.Range("C82").Value = "Sim, existe compartilhamento dentro da nuvem"
$ws.Range("C83").Value = "Tem acesso, porém é efetuado apenas mediante aviso e autorização"
$ws.Range("C84").Value = "Virtualização para cada cliente e ambiente"

# B80 label text reverts to the original question wording
$ws.Range("B80").Value = "Todos os módulos da solução são atualizados sempre na mesma versão?"
